$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header B1 from "y" to "dy0"
$ws.Range("B1").Value = "dy0"

# 2. Rename J4 from "h" to "xi"
$ws.Range("J4").Value = "xi"

# 3. Add new row 8 entries
$ws.Range("J8").Value = "t_2"
$ws.Range("K8").Formula = "=(K2-A5)/K4"

# 4. Apply centered alignment (horizontal + vertical) to the whole B1:H8 block.
#    Set it on a single cell first, then paste the format onto the whole range
#    so only one combined style gets created (instead of two separate ones).
$fmtCell = $ws.Range("B1")
$fmtCell.HorizontalAlignment = -4108   # xlCenter
$fmtCell.VerticalAlignment = -4108     # xlCenter
$blockRng = $ws.Range("B1:H8")
$fmtCell.Copy()
$blockRng.PasteSpecial(-4122)          # xlPasteFormats

# 5. Highlight the diagonal cells with a yellow fill, on top of the centered
#    alignment, again combining both format changes into a single style via copy/paste.
$hlCell = $ws.Range("G2")
$hlCell.HorizontalAlignment = -4108
$hlCell.VerticalAlignment = -4108
$hlCell.Interior.Color = 65535          # yellow
$hlRng = $ws.Range("G2:H2,E3:F3,C4:D4,B5")
$hlCell.Copy()
$hlRng.PasteSpecial(-4122)

# 6. Leave the final selection on B1:H8, matching the workbook's saved selection state
$excel.CutCopyMode = $false
$blockRng.Select
